# ChampTable.docx edit script
# Implements:
#  - "Clients" -> "Utilisateur" (section title)
#  - "Username" -> "N Utilisateur" (degree sign) in the Utilisateur block
#  - insert "Administrateur ?" line after "Prenom" in the Utilisateur person block
#  - Ticket block: add "N Client" line, drop "Username" line, annotate IP / Nom
#    machine with " (cache)", add "Commentaire" line before "Resolution"
#  - append two new paragraphs ("Administrateur" title + "N utilisateur /
#    Niveau de droit" body)

$d = $word.ActiveDocument

# --- "Clients" title -> "Utilisateur" ---------------------------------
$p = $d.Paragraphs(5).Range
$p.Find.Execute("Clients", $true, $false, $false, $false, $false, $true, 1, $false, "Utilisateur", 2) | Out-Null

# --- "Username" -> "N(deg) Utilisateur" in the username/nom paragraph --
$p = $d.Paragraphs(6).Range
$p.Find.Execute("Username", $true, $false, $false, $false, $false, $true, 1, $false, "N° Utilisateur", 2) | Out-Null

# --- insert "Administrateur ?" line after "Prenom" in the person block -
$p = $d.Paragraphs(7).Range
$p.Find.Execute("Prénom", $true, $false, $false, $false, $false, $true, 1, $false, "Prénom`vAdministrateur ?", 2) | Out-Null

# --- Ticket fields paragraph -------------------------------------------
$ticketIndex = 9

# a) "N de ticket" -> add "N Client" line right after it
$p = $d.Paragraphs($ticketIndex).Range
$p.Find.Execute("N° de ticket", $true, $false, $false, $false, $false, $true, 1, $false, "N° de ticket`vN° Client", 2) | Out-Null

# b) drop the "Username" line entirely (text + its trailing line break)
$p = $d.Paragraphs($ticketIndex).Range
$p.Find.Execute("Username`v", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# c) "IP" -> "IP (caché)"
$p = $d.Paragraphs($ticketIndex).Range
$p.Find.Execute("IP`v", $true, $false, $false, $false, $false, $true, 1, $false, "IP (caché)`v", 2) | Out-Null

# d) "Nom machine" -> "Nom machine (caché)"
$p = $d.Paragraphs($ticketIndex).Range
$p.Find.Execute("Nom machine", $true, $false, $false, $false, $false, $true, 1, $false, "Nom machine (caché)", 2) | Out-Null

# e) insert "Commentaire" line after "Description", before "Résolution"
$p = $d.Paragraphs($ticketIndex).Range
$p.Find.Execute("Description", $true, $false, $false, $false, $false, $true, 1, $false, "Description`vCommentaire", 2) | Out-Null

# f) relocate the "_GoBack" bookmark so it sits right after "Commentaire"
#    (just like it originally sat right after "Nom machine", before "Statut")
$p = $d.Paragraphs($ticketIndex).Range
$p.Find.Execute("Commentaire") | Out-Null
$p.Collapse(0)
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$d.Bookmarks.Add("_GoBack", $p)

# --- append the new "Administrateur" title paragraph --------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newTitle = $d.Paragraphs($d.Paragraphs.Count)
$newTitle.Range.Text = "Administrateur"
$newTitle.Style = "Titre1"

# --- append the new "N utilisateur / Niveau de droit" body paragraph ----
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newBody = $d.Paragraphs($d.Paragraphs.Count)
$newBody.Style = "Normal"
$newBody.Range.Text = "N° utilisateur`vNiveau de droit"
